$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds plain-text figures (e.g. "28.075.39", "0.00001063") that are
# not valid numeric literals, so the whole data range is kept as Text-formatted before
# writing the refreshed values -- otherwise Excel would silently reinterpret look-alike
# numbers (e.g. "1.004") as numbers and drop formatting such as trailing zeros.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.028.35'
$ws.Range("E2").Value = '  -2.79%  '
$ws.Range("D3").Value = '1.896.15'
$ws.Range("E3").Value = '  -3.51%  '
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.98%  '
$ws.Range("D5").Value = '326.56'
$ws.Range("E5").Value = '  +0.90%  '
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  -0.85%  '
$ws.Range("D7").Value = '0.4608'
$ws.Range("E7").Value = '  -3.33%  '
$ws.Range("D8").Value = '0.3951'
$ws.Range("E8").Value = '  -2.28%  '
$ws.Range("D9").Value = '51.73'
$ws.Range("E9").Value = '  -4.23%  '
$ws.Range("D10").Value = '0.08275'
$ws.Range("E10").Value = '  -3.33%  '
$ws.Range("D11").Value = '1.038'
$ws.Range("E11").Value = '  -2.09%  '
$ws.Range("D12").Value = '21.72'
$ws.Range("E12").Value = '  -3.44%  '
$ws.Range("D13").Value = '1.874.59'
$ws.Range("E13").Value = '  -6.41%  '
$ws.Range("D14").Value = '7.332'
$ws.Range("E14").Value = '  -4.12%  '
$ws.Range("D15").Value = '5.988'
$ws.Range("E15").Value = '  -4.10%  '
$ws.Range("D16").Value = '1.006'
$ws.Range("E16").Value = '  -0.83%  '
$ws.Range("D17").Value = '89.37'
$ws.Range("E17").Value = '  -0.84%  '
$ws.Range("D18").Value = '0.00001062'
$ws.Range("E18").Value = '  -0.76%  '
$ws.Range("D19").Value = '0.06588'
$ws.Range("E19").Value = '  -0.60%  '
$ws.Range("D20").Value = '17.68'
$ws.Range("E20").Value = '  -5.37%  '
$ws.Range("E21").Value = '  -1.28%  '
$ws.Range("D22").Value = '5.672'
$ws.Range("E22").Value = '  -1.98%  '
$ws.Range("D23").Value = '28.013.53'
$ws.Range("E23").Value = '  -2.89%  '
$ws.Range("D24").Value = '11.10'
$ws.Range("E24").Value = '  -3.97%  '
$ws.Range("D25").Value = '2.305'
$ws.Range("E25").Value = '  +0.44%  '
$ws.Range("D26").Value = '2.105.28'
$ws.Range("E26").Value = '  -6.12%  '
$ws.Range("D27").Value = '154.18'
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("D28").Value = '19.95'
$ws.Range("E28").Value = '  -1.49%  '
$ws.Range("D29").Value = '2.108'
$ws.Range("E29").Value = '  -2.44%  '
$ws.Range("D30").Value = '5.681'
$ws.Range("E30").Value = '  -4.86%  '
$ws.Range("D31").Value = '124.65'
$ws.Range("E31").Value = '  +0.28%  '
$ws.Range("D32").Value = '0.09570'
$ws.Range("E32").Value = '  -0.67%  '
$ws.Range("D33").Value = '0.9630'
$ws.Range("E33").Value = '  -4.35%  '
$ws.Range("D34").Value = '1.458'
$ws.Range("E34").Value = '  -0.17%  '
$ws.Range("D35").Value = '3.627'
$ws.Range("E35").Value = '  -1.82%  '
$ws.Range("D36").Value = '5.468'
$ws.Range("E36").Value = '  -3.72%  '
$ws.Range("D37").Value = '1.259'
$ws.Range("E37").Value = '  -0.47%  '
$ws.Range("D38").Value = '0.02287'
$ws.Range("E38").Value = '  -2.45%  '
$ws.Range("D39").Value = '8.654'
$ws.Range("E39").Value = '  -1.28%  '
$ws.Range("D40").Value = '0.06110'
$ws.Range("E40").Value = '  -1.75%  '
$ws.Range("D41").Value = '0.6102'
$ws.Range("E41").Value = '  -2.37%  '
$ws.Range("E42").Value = '  -1.08%  '
$ws.Range("D43").Value = '10.76'
$ws.Range("E43").Value = '  -3.19%  '
$ws.Range("D44").Value = '0.1894'
$ws.Range("E44").Value = '  -1.53%  '
$ws.Range("D45").Value = '1.304'
$ws.Range("E45").Value = '  -2.53%  '
$ws.Range("D48").Value = '1.997'
$ws.Range("E48").Value = '  -4.18%  '
$ws.Range("D49").Value = '3.430'
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("D50").Value = '0.06892'
$ws.Range("E50").Value = '  +0.53%  '
$ws.Range("D51").Value = '110.53'
$ws.Range("E51").Value = '  -0.32%  '

# Rows 46 and 47 swapped places (coin identity + its associated price/volume moved rows),
# and the price/volume figures for both coins were also refreshed.
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '12.81'
$ws.Range("E46").Value = '  -1.04%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '0.5826'
$ws.Range("E47").Value = '  -2.34%  '
